$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full recomputed LR-pair table (Sending cluster x Target cluster, all
# combinations of ECs/FAPs/M2/sCs) after following Dr Hou's advice.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the expression / specificity metrics.
$data = @(
  @("ECs","Hspg2","Itga2","ECs",3,1,135.7046583333333,407.113975,0.2901853119378819,0.2901853119378819,2,0.6666666666666666,2.672731,8.018193,0.5408378022089502,0.5408378022089503,362.7020471719084,3264.318424547175,0.1569431863418027,0.1569431863418027),
  @("ECs","Hspg2","Itga2","FAPs",3,1,135.7046583333333,407.113975,0.2901853119378819,0.2901853119378819,3,1,1.627877666666667,4.883633,0.3294075533620981,0.3294075533620982,220.9105825634639,1988.195243071175,0.09558923362707493,0.09558923362707496),
  @("ECs","Hspg2","Itga2","M2",3,1,135.7046583333333,407.113975,0.2901853119378819,0.2901853119378819,1,0.3333333333333333,0.01852966666666667,0.055589,0.003749552123152102,0.003749552123152104,2.514562084030556,22.631058756275,0.00108806495248424,0.001088064952484241),
  @("ECs","Hspg2","Itga2","sCs",3,1,135.7046583333333,407.113975,0.2901853119378819,0.2901853119378819,3,1,0.6226963333333334,1.868089,0.1260050923057995,0.1260050923057995,84.50279316041946,760.5251384437751,0.03656482701652004,0.03656482701652004),
  @("FAPs","Hspg2","Itga2","ECs",3,1,226.082006,678.246018,0.4834445496594812,0.4834445496594812,2,0.6666666666666666,2.672731,8.018193,0.5408378022089502,0.5408378022089503,604.2563859783861,5438.307473805475,0.2614650877277295,0.2614650877277296),
  @("FAPs","Hspg2","Itga2","FAPs",3,1,226.082006,678.246018,0.4834445496594812,0.4834445496594812,3,1,1.627877666666667,4.883633,0.3294075533620981,0.3294075533620982,368.0338484025993,3312.304635623394,0.1592502862895711,0.1592502862895711),
  @("FAPs","Hspg2","Itga2","M2",3,1,226.082006,678.246018,0.4834445496594812,0.4834445496594812,1,0.3333333333333333,0.01852966666666667,0.055589,0.003749552123152102,0.003749552123152104,4.189224210511333,37.70301789460201,0.00181270053760202,0.001812700537602021),
  @("FAPs","Hspg2","Itga2","sCs",3,1,226.082006,678.246018,0.4834445496594812,0.4834445496594812,3,1,0.6226963333333334,1.868089,0.1260050923057995,0.1260050923057995,140.7804361688447,1267.023925519602,0.06091647510457861,0.06091647510457862),
  @("M2","Hspg2","Itga2","ECs",3,1,0.1600446666666666,0.4801339999999999,0.0003422329939962955,0.0003422329939962955,2,0.6666666666666666,2.672731,8.018193,0.5408378022089502,0.5408378022089503,0.4277563419846666,3.849807077862,0.0001850925403163453,0.0001850925403163453),
  @("M2","Hspg2","Itga2","FAPs",3,1,0.1600446666666666,0.4801339999999999,0.0003422329939962955,0.0003422329939962955,3,1,1.627877666666667,4.883633,0.3294075533620981,0.3294075533620982,0.2605331385357777,2.344798246822,0.0001127341332321053,0.0001127341332321053),
  @("M2","Hspg2","Itga2","M2",3,1,0.1600446666666666,0.4801339999999999,0.0003422329939962955,0.0003422329939962955,1,0.3333333333333333,0.01852966666666667,0.055589,0.003749552123152102,0.003749552123152104,0.002965574325111111,0.026690168926,0.00000128322044925151,0.000001283220449251511),
  @("M2","Hspg2","Itga2","sCs",3,1,0.1600446666666666,0.4801339999999999,0.0003422329939962955,0.0003422329939962955,3,1,0.6226963333333334,1.868089,0.1260050923057995,0.1260050923057995,0.09965922710288888,0.8969330439259999,0.00004312309999859334,0.00004312309999859335),
  @("sCs","Hspg2","Itga2","ECs",3,1,105.7015583333333,317.104675,0.2260279054086406,0.2260279054086406,2,0.6666666666666666,2.672731,8.018193,0.5408378022089502,0.5408378022089503,282.5118317058083,2542.606485352275,0.1222444355991017,0.1222444355991017),
  @("sCs","Hspg2","Itga2","FAPs",3,1,105.7015583333333,317.104675,0.2260279054086406,0.2260279054086406,3,1,1.627877666666667,4.883633,0.3294075533620981,0.3294075533620982,172.0692061426972,1548.622855284275,0.07445529931222003,0.07445529931222006),
  @("sCs","Hspg2","Itga2","M2",3,1,105.7015583333333,317.104675,0.2260279054086406,0.2260279054086406,1,0.3333333333333333,0.01852966666666667,0.055589,0.003749552123152102,0.003749552123152104,1.958614642063889,17.627531778575,0.000847503412616591,0.0008475034126165912),
  @("sCs","Hspg2","Itga2","sCs",3,1,105.7015583333333,317.104675,0.2260279054086406,0.2260279054086406,3,1,0.6226963333333334,1.868089,0.1260050923057995,0.1260050923057995,65.81997280178611,592.379755216075,0.02848066708470228,0.02848066708470228)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
